$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.55"
$ws.Range("E2").Value = "'-5.43%"
$ws.Range("E3").Value = "'-1.89%"
$ws.Range("D4").Value = "'5.036"
$ws.Range("E4").Value = "'-3.03%"
$ws.Range("D5").Value = "'0.07359"
$ws.Range("E5").Value = "'-4.17%"
$ws.Range("D6").Value = "'4.295"
$ws.Range("E6").Value = "'-0.27%"
$ws.Range("D7").Value = "'1.553"
$ws.Range("E7").Value = "'-7.50%"
$ws.Range("D8").Value = "'0.9176"
$ws.Range("E8").Value = "'0.30%"
$ws.Range("D9").Value = "'0.1189"
$ws.Range("E9").Value = "'-3.96%"
$ws.Range("D10").Value = "'0.1736"
$ws.Range("E10").Value = "'-4.97%"
$ws.Range("D11").Value = "'0.08755"
$ws.Range("E11").Value = "'-4.09%"
$ws.Range("D12").Value = "'0.04200"
$ws.Range("E12").Value = "'0.37%"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("D14").Value = "'0.001272"
$ws.Range("E14").Value = "'-2.02%"
$ws.Range("D15").Value = "'0.005758"
$ws.Range("E15").Value = "'-0.20%"
$ws.Range("D16").Value = "'3.409"
$ws.Range("E16").Value = "'1.90%"
$ws.Range("E17").Value = "'-1.07%"
$ws.Range("D18").Value = "'0.3297"
$ws.Range("E18").Value = "'-0.67%"
$ws.Range("D19").Value = "'7.604"
$ws.Range("E19").Value = "'2.61%"
$ws.Range("D20").Value = "'0.1353"
$ws.Range("E20").Value = "'-1.71%"
$ws.Range("E21").Value = "'1.25%"
$ws.Range("D22").Value = "'0.03848"
$ws.Range("E22").Value = "'-4.23%"
$ws.Range("D23").Value = "'0.001282"
$ws.Range("E23").Value = "'0.98%"
$ws.Range("D24").Value = "'0.003891"
$ws.Range("E24").Value = "'-4.81%"
$ws.Range("D25").Value = "'0.0001283"
$ws.Range("E25").Value = "'-1.64%"
$ws.Range("D38").Value = "'0.02341"
$ws.Range("E38").Value = "'-7.32%"
$ws.Range("D39").Value = "'0.05036"
$ws.Range("E39").Value = "'-5.22%"
$ws.Range("D40").Value = "'0.007692"
$ws.Range("E40").Value = "'-1.90%"
$ws.Range("E41").Value = "'166.91%"
$ws.Range("D42").Value = "'0.1271"
$ws.Range("E42").Value = "'-2.76%"
$ws.Range("D43").Value = "'0.007374"
$ws.Range("D44").Value = "'0.007677"
$ws.Range("E44").Value = "'3.56%"
$ws.Range("D45").Value = "'0.3178"
$ws.Range("E45").Value = "'3.65%"
$ws.Range("D46").Value = "'0.00006520"
$ws.Range("E46").Value = "'-4.16%"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("E48").Value = "'7.85%"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("E51").Value = "'-0.13%"
